# Daily attendance processing - 2026-01-27 01:47:03
# Reorders the "Recorded By" (column G) values for specific recorder-list
# combinations: swap the first two comma-separated entries, leaving any
# trailing entries untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Exact "before" values that need their first two entries swapped, mapped
# to the resulting "after" value.
$swapMap = @{
    "System, dnasr281@gmail.com"           = "dnasr281@gmail.com, System"
    "backup@backdoor.com, system, System"  = "system, backup@backdoor.com, System"
    "admin@admin.com, dnasr281@gmail.com"  = "dnasr281@gmail.com, admin@admin.com"
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    if ($swapMap.ContainsKey($val)) {
        $cell.Value = $swapMap[$val]
    }
}
